$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: was a numeric value, becomes blank (empty text) like the other
# unused grid cells on this sheet. Using a leading apostrophe forces the
# literal to commit as an empty text value rather than clearing the cell
# to a blank/number; copying C3's style afterwards drops the resulting
# quote-prefix formatting so D3 ends up on the default (unstyled) cell
# format, matching its neighbours.
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = $ws.Range("C3").Style

# C4 / C5: numeric values updated with corrected figures.
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 1922.615322677455

# Row 7 is renamed from "Other" to "Biogas" and its value corrected.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 610.5425936133935

# A new row 8 is inserted carrying the "Other" label (the row that used
# to live at r=7) together with its own figure. Copy the formatting of
# row 7's label cell first so the new label cell picks up the same
# (bold/bordered/centered) style index, then overwrite the value.
$ws.Range("A7").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "Other"

# B8 / C8 stay blank (empty text), same convention as every other unused
# cell in this grid (e.g. B7, C7).
$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = $ws.Range("B7").Style
$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = $ws.Range("C7").Style

$ws.Range("D8").Value = 1713.093202559996
